$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor recalculated timestamp value on the previous last row (A9)
$ws.Cells.Item(9, 1).Value = 45865.37525619213

# New row of sensor data appended by the scheduled task (row 10)
$ws.Cells.Item(10, 1).Value = 45865.41694434191
$ws.Cells.Item(10, 2).Value = 2025
$ws.Cells.Item(10, 3).Value = 30
$ws.Cells.Item(10, 4).Value = 14.72
$ws.Cells.Item(10, 5).Value = 86.02
$ws.Cells.Item(10, 6).Value = 240.82
$ws.Cells.Item(10, 7).Value = 9.460000000000001
$ws.Cells.Item(10, 8).Value = "ESE"
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = "10:00:23"

# Match the date/time number format used by the other rows in column A
$ws.Cells.Item(10, 1).NumberFormat = $ws.Cells.Item(9, 1).NumberFormat
